# Revert 'cards' to commit 95cda46ab8 (Jun 25)
#
# Original layout: A=pst_code, B=pst_labe, C=updated_at
# Target layout:   A=_airbyte_ab_id, B=_airbyte_emitted_at, C=pst_code, D=pst_labe,
#                  E=_airbyte_additional_properties, F=source_file_path, G=updated_at

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns in front (A,B) for _airbyte_ab_id / _airbyte_emitted_at.
# This pushes pst_code->C, pst_labe->D, updated_at->E (shifting content+formats along).
$ws.Range("A:B").Insert()

# Insert 2 more new columns (E,F) for _airbyte_additional_properties / source_file_path,
# right before the (now shifted) updated_at column, pushing updated_at->G.
$ws.Range("E:F").Insert()

# ---- Header row ----
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("C1").Value = "pst_code"
$ws.Range("D1").Value = "pst_labe"
$ws.Range("E1").Value = "_airbyte_additional_properties"
$ws.Range("F1").Value = "source_file_path"
$ws.Range("G1").Value = "updated_at"

# New header cells come in unformatted - copy the bold/border/center header style
# from a cell that already carries it (C1, the original pst_code header) onto them.
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)

# ---- Row 2 ----
$ws.Range("A2").Value = "431a9c04-da43-4101-8bbb-2cded7053cfb"
$ws.Range("B2").Value = 45510.3079196875
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/POS_TERMINAL_STATUS/2024_08_06_1722929004063_0.parquet"
$ws.Range("G2").Value = 45511.29527901154

# ---- Row 3 ----
$ws.Range("A3").Value = "0e96ddbe-c889-48cb-bd2b-497b27b4d8d2"
$ws.Range("B3").Value = 45510.3079196875
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/POS_TERMINAL_STATUS/2024_08_06_1722929004063_0.parquet"
$ws.Range("G3").Value = 45511.29527901154

# The new _airbyte_emitted_at column (B) should carry the same date/time number
# format as the updated_at column (G) - copy that format across.
$ws.Range("G2").Copy()
$ws.Range("B2:B3").PasteSpecial(-4122)
